{"js": "// Fix the copy and layout of About page:\n// The \"Skills\" paragraph starts with a bold \"JavaScript\" run. We need to:\n//   1. Turn that run's text into \"JavaScript \" (trailing space added, still bold).\n//   2. Insert a brand-new, NOT-bold run right after it containing \"(ES6)\".\n// Net visible result: \"JavaScript (ES6), TypeScript, React, HTML5, ...\"\n//\n// There are other \"JavaScript\" occurrences in the document (an intro sentence,\n// and \"JavaScript animation\" later in the same Skills paragraph), so we anchor\n// on the unique, longer substring \"JavaScript, TypeScript, React, HTML5\" that\n// only exists at the very start of the Skills list before narrowing down to\n// the \"JavaScript\" word itself.\n\nconst body = context.document.body;\n\nconst anchorResults = body.search(\"JavaScript, TypeScript, React, HTML5\", {\n  matchCase: true\n});\nanchorResults.load(\"items\");\nawait context.sync();\n\nif (anchorResults.items.length === 0) {\n  throw new Error(\"Could not find the Skills paragraph anchor text.\");\n}\n\nconst anchorRange = anchorResults.items[0];\nconst jsResults = anchorRange.search(\"JavaScript\", { matchCase: true });\njsResults.load(\"items\");\nawait context.sync();\n\nif (jsResults.items.length === 0) {\n  throw new Error(\"Could not find the 'JavaScript' run inside the Skills paragraph.\");\n}\n\n// First hit within the anchor range is the standalone, bold \"JavaScript\" term.\nconst jsRange = jsResults.items[0];\n\n// 1) Append a trailing space \u2014 lands in the same (bold) run as \"JavaScript\".\nconst spaceRange = jsRange.insertText(\" \", \"End\");\n\n// 2) Insert \"(ES6)\" right after that space as its own run, then turn bold off.\nconst es6Range = spaceRange.insertText(\"(ES6)\", \"End\");\nes6Range.font.bold = false;\n\nawait context.sync();\n", "ps1": "# Fix the copy and layout of About page:\n# The \"Skills\" paragraph starts with a bold \"JavaScript\" run. We need to:\n#   1. Turn that run's text into \"JavaScript \" (trailing space added, still bold).\n#   2. Insert a brand-new, NOT-bold run right after it containing \"(ES6)\".\n# Net visible result: \"JavaScript (ES6), TypeScript, React, HTML5, ...\"\n#\n# There are other \"JavaScript\" occurrences in the document (an intro sentence,\n# and \"JavaScript animation\" later in the same Skills paragraph), so we first\n# anchor on the unique, longer substring \"JavaScript, TypeScript, React, HTML5\"\n# that only exists at the very start of the Skills list, then scope a second\n# Find to just that anchor range to land on the standalone \"JavaScript\" term.\n\n$d = $word.ActiveDocument\n\n# Step 1: locate the unique anchor text that starts the Skills paragraph.\n$anchor = $d.Content\n$anchor.Find.Text = \"JavaScript, TypeScript, React, HTML5\"\n$anchor.Find.MatchCase = $true\n$anchor.Find.Execute() | Out-Null\n\n# Step 2: find \"JavaScript\" scoped within that anchor range only (so the\n# bold, standalone term is matched, not the \"JavaScript animation\" run).\n$jsRng = $d.Range($anchor.Start, $anchor.End)\n$jsRng.Find.Text = \"JavaScript\"\n$jsRng.Find.MatchCase = $true\n$jsRng.Find.Execute() | Out-Null\n\n# 1) Append a trailing space \u2014 lands in the same (bold) run as \"JavaScript\".\n$jsRng.InsertAfter(\" \")\n\n# 2) Insert \"(ES6)\" right after that space as its own run, then turn bold off.\n$esRng = $d.Range($jsRng.End, $jsRng.End)\n$esRng.InsertAfter(\"(ES6)\")\n$esRng.Bold = 0\n"}
